# Ignorar carpeta de sesión en Git
# Append the newest batch of WhatsApp-bot conversation rows (ProspectoBot)
# to the "Mensajes" sheet: Fecha, Hora, Numero, Nombre, Mensaje.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mensajes")

$rows = @(
    @("8/5/2025", "12:47:59 p. m.", "51955091726@c.us", "~Shukita", "info sap sd"),
    @("8/5/2025", "2:17:57 p. m.",  "51955091726@c.us", "~Shukita", "info de sql server"),
    @("8/5/2025", "2:18:30 p. m.",  "51955091726@c.us", "~Shukita", "miau y si deseo la inscripcion"),
    @("8/5/2025", "2:18:49 p. m.",  "51955091726@c.us", "~Shukita", "sap fi"),
    @("8/5/2025", "2:18:58 p. m.",  "51955091726@c.us", "~Shukita", "Excel basico")
)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$startRow = $lastRow + 1

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Column A holds a date-looking string ("8/5/2025"). Excel's normal
    # value-assignment auto-converts such text into a real date serial,
    # so force the cell to Text first, assign the literal string, then
    # restore the cell's style so it still matches the rest of the sheet.
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $data[0]
    $dateCell.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
}
